$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2984.2273
$ws.Range("I40").Value = 2825.125
$ws.Range("J40").Value = 3075.1428
$ws.Range("K40").Value = 2825.125
$ws.Range("L40").Value = 3075.1428
$ws.Range("M40").Value = -2650.125
$ws.Range("N40").Value = -3425.1428
$ws.Range("H51").Value = 4585.7144
$ws.Range("I51").Value = 3150
$ws.Range("K51").Value = 3150
$ws.Range("M51").Value = -2666
$ws.Range("H64").Value = 145314.28
$ws.Range("J64").Value = 2840
$ws.Range("L64").Value = 2840
$ws.Range("N64").Value = -3336
$ws.Range("H67").Value = 145314.28
$ws.Range("J67").Value = 2840
$ws.Range("L67").Value = 2840
$ws.Range("N67").Value = -4556
$ws.Range("H70").Value = 998.3333
$ws.Range("I70").Value = 998.3333
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2994.9999
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2724.9999
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 998.3333
$ws.Range("I73").Value = 998.3333
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2994.9999
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -2058.9999
$ws.Range("N73").ClearContents()
$ws.Range("H109").Value = 35758.5
$ws.Range("J109").Value = 35758.5
$ws.Range("L109").Value = 35758.5
$ws.Range("N109").Value = -38532.5
$ws.Range("H111").Value = 1582.25
$ws.Range("I111").Value = 1582.25
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 4746.75
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -1679.75
$ws.Range("N111").ClearContents()
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H129").Value = 1631.2441
$ws.Range("I129").Value = 1459
$ws.Range("J129").Value = 1664.7361
$ws.Range("K129").Value = 4377
$ws.Range("L129").Value = 4994.2083
$ws.Range("M129").Value = 623
$ws.Range("N129").Value = -14994.2083
$ws.Range("H138").Value = 1975.2872
$ws.Range("I138").Value = 1472.5122
$ws.Range("J138").Value = 2364.2263
$ws.Range("K138").Value = 4417.536599999999
$ws.Range("L138").Value = 7092.678899999999
$ws.Range("M138").Value = 722.4634000000005
$ws.Range("N138").Value = -17372.6789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19404.285
$ws.Range("I32").Value = 19970.865
$ws.Range("K32").Value = 19970.865
$ws.Range("M32").Value = -19683.865
$ws.Range("H38").Value = 10000
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10934
$ws.Range("H107").Value = 37369.75
$ws.Range("J107").Value = 37369.75
$ws.Range("L107").Value = 37369.75
$ws.Range("N107").Value = -45049.75
$ws.Range("H123").Value = 51429
$ws.Range("J123").Value = 51429
$ws.Range("L123").Value = 51429
$ws.Range("N123").Value = -61229
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 3845.6785
$ws.Range("I80").Value = 14380
$ws.Range("J80").Value = 334.2381
$ws.Range("K80").Value = 14380
$ws.Range("L80").Value = 334.2381
$ws.Range("M80").Value = -13382
$ws.Range("N80").Value = -2330.2381
$ws.Range("H83").Value = 3845.6785
$ws.Range("I83").Value = 14380
$ws.Range("J83").Value = 334.2381
$ws.Range("K83").Value = 71900
$ws.Range("L83").Value = 1671.1905
$ws.Range("M83").Value = -66908
$ws.Range("N83").Value = -11655.1905
$ws.Range("H117").Value = 49997.332
$ws.Range("J117").Value = 49997.332
$ws.Range("L117").Value = 49997.332
$ws.Range("N117").Value = -59175.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3250
$ws.Range("I62").Value = 3250
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3250
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2626
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3250
$ws.Range("I65").Value = 3250
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16250
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -13130
$ws.Range("N65").ClearContents()
$ws.Range("H99").Value = 1923.1111
$ws.Range("I99").Value = 1630.2858
$ws.Range("K99").Value = 1630.2858
$ws.Range("M99").Value = -132.2858000000001
$ws.Range("H112").Value = 38948.2
$ws.Range("J112").Value = 38948.2
$ws.Range("L112").Value = 38948.2
$ws.Range("N112").Value = -41902.2
$ws.Range("H116").Value = 47814.332
$ws.Range("J116").Value = 47814.332
$ws.Range("L116").Value = 47814.332
$ws.Range("N116").Value = -56992.332
$ws.Range("H126").Value = 1923.1111
$ws.Range("I126").Value = 1630.2858
$ws.Range("K126").Value = 4890.857400000001
$ws.Range("M126").Value = -2420.857400000001
$ws.Range("H130").Value = 37429.43
$ws.Range("J130").Value = 37429.43
$ws.Range("L130").Value = 37429.43
$ws.Range("N130").Value = -47469.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1415.1714
$ws.Range("I5").Value = 612.06665
$ws.Range("J5").Value = 2017.5
$ws.Range("K5").Value = 1836.19995
$ws.Range("L5").Value = 6052.5
$ws.Range("M5").Value = -1724.19995
$ws.Range("N5").Value = -6276.5
$ws.Range("H23").Value = 499.36365
$ws.Range("I23").Value = 390.33334
$ws.Range("J23").Value = 540.25
$ws.Range("K23").Value = 1171.00002
$ws.Range("L23").Value = 1620.75
$ws.Range("M23").Value = -936.0000199999999
$ws.Range("N23").Value = -2090.75
$ws.Range("H122").Value = 2716.7646
$ws.Range("I122").Value = 662.5227
$ws.Range("J122").Value = 15629.143
$ws.Range("K122").Value = 5962.704299999999
$ws.Range("L122").Value = 140662.287
$ws.Range("M122").Value = -3512.704299999999
$ws.Range("N122").Value = -145562.287
$ws.Range("H135").Value = 1415.1714
$ws.Range("I135").Value = 612.06665
$ws.Range("J135").Value = 2017.5
$ws.Range("K135").Value = 5508.59985
$ws.Range("L135").Value = 18157.5
$ws.Range("M135").Value = -2973.59985
$ws.Range("N135").Value = -23227.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 49434
$ws.Range("J116").Value = 49434
$ws.Range("L116").Value = 49434
$ws.Range("N116").Value = -58612
$ws.Range("H119").Value = 48761
$ws.Range("J119").Value = 48761
$ws.Range("L119").Value = 48761
$ws.Range("N119").Value = -58437

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 4020
$ws.Range("I30").Value = 4020
$ws.Range("K30").Value = 4020
$ws.Range("M30").Value = -3912
$ws.Range("H68").Value = 2992.3076
$ws.Range("I68").Value = 2725
$ws.Range("J68").Value = 3420
$ws.Range("K68").Value = 2725
$ws.Range("L68").Value = 3420
$ws.Range("M68").Value = -1976
$ws.Range("N68").Value = -4918
$ws.Range("H71").Value = 2992.3076
$ws.Range("I71").Value = 2725
$ws.Range("J71").Value = 3420
$ws.Range("K71").Value = 13625
$ws.Range("L71").Value = 17100
$ws.Range("M71").Value = -9881
$ws.Range("N71").Value = -24588
$ws.Range("H114").Value = 38344
$ws.Range("J114").Value = 38344
$ws.Range("L114").Value = 38344
$ws.Range("N114").Value = -47022
$ws.Range("H116").Value = 50330.5
$ws.Range("J116").Value = 50330.5
$ws.Range("L116").Value = 50330.5
$ws.Range("N116").Value = -59508.5
$ws.Range("H117").Value = 41786
$ws.Range("J117").Value = 41786
$ws.Range("L117").Value = 41786
$ws.Range("N117").Value = -50964

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45775.5
$ws.Range("J16").Value = 45775.5
$ws.Range("L16").Value = 45775.5
$ws.Range("N16").Value = -46359.5
$ws.Range("H38").Value = 4124.75
$ws.Range("J38").Value = 4124.75
$ws.Range("L38").Value = 4124.75
$ws.Range("N38").Value = -5070.75
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H81").Value = 2176.6667
$ws.Range("I81").Value = 2012
$ws.Range("K81").Value = 4024
$ws.Range("M81").Value = -2963
$ws.Range("H84").Value = 2176.6667
$ws.Range("I84").Value = 2012
$ws.Range("K84").Value = 20120
$ws.Range("M84").Value = -14816
$ws.Range("H117").Value = 47699.668
$ws.Range("J117").Value = 47699.668
$ws.Range("L117").Value = 47699.668
$ws.Range("N117").Value = -56877.668
$ws.Range("H119").Value = 49690
$ws.Range("J119").Value = 49690
$ws.Range("L119").Value = 49690
$ws.Range("N119").Value = -59366
